{"js": "const replacements = [\n  [\"2025-11-29 Saturday\", \"2025-11-30 Sunday\"],\n  [\"291\u00d75=\", \"464\u00d76=\"],\n  [\"263\u00d76=\", \"656\u00d79=\"],\n  [\"913\u00d76=\", \"789\u00d79=\"],\n  [\"173\u00d73=\", \"802\u00d79=\"],\n  [\"417\u00d77=\", \"550\u00d75=\"],\n  [\"754\u00d75=\", \"629\u00d76=\"],\n  [\"446\u00d76=\", \"798\u00d75=\"],\n  [\"939\u00d73=\", \"502\u00d79=\"],\n  [\"562\u00d75=\", \"927\u00d73=\"],\n  [\"452\u00d77=\", \"251\u00d79=\"],\n  [\"456\u00d75=\", \"172\u00d73=\"],\n  [\"778\u00d72=\", \"787\u00d73=\"],\n  [\"301\u00d73=\", \"204\u00d74=\"],\n  [\"670\u00d72=\", \"181\u00d77=\"],\n  [\"626\u00d73=\", \"492\u00d73=\"],\n  [\"102\u00d78=\", \"659\u00d77=\"],\n  [\"664\u00d76=\", \"511\u00d78=\"],\n  [\"223\u00d73=\", \"399\u00d74=\"],\n  [\"626\u00d75=\", \"874\u00d72=\"],\n  [\"428\u00d72=\", \"766\u00d78=\"],\n  [\"954\u00d77=\", \"336\u00d72=\"],\n  [\"297\u00d72=\", \"149\u00d75=\"],\n  [\"464\u00d77=\", \"573\u00d74=\"],\n  [\"679\u00d75=\", \"221\u00d76=\"],\n  [\"913\u00d75=\", \"322\u00d74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Find = \"2025-11-29 Saturday\"; Replace = \"2025-11-30 Sunday\"},\n    @{Find = \"291\u00d75=\"; Replace = \"464\u00d76=\"},\n    @{Find = \"263\u00d76=\"; Replace = \"656\u00d79=\"},\n    @{Find = \"913\u00d76=\"; Replace = \"789\u00d79=\"},\n    @{Find = \"173\u00d73=\"; Replace = \"802\u00d79=\"},\n    @{Find = \"417\u00d77=\"; Replace = \"550\u00d75=\"},\n    @{Find = \"754\u00d75=\"; Replace = \"629\u00d76=\"},\n    @{Find = \"446\u00d76=\"; Replace = \"798\u00d75=\"},\n    @{Find = \"939\u00d73=\"; Replace = \"502\u00d79=\"},\n    @{Find = \"562\u00d75=\"; Replace = \"927\u00d73=\"},\n    @{Find = \"452\u00d77=\"; Replace = \"251\u00d79=\"},\n    @{Find = \"456\u00d75=\"; Replace = \"172\u00d73=\"},\n    @{Find = \"778\u00d72=\"; Replace = \"787\u00d73=\"},\n    @{Find = \"301\u00d73=\"; Replace = \"204\u00d74=\"},\n    @{Find = \"670\u00d72=\"; Replace = \"181\u00d77=\"},\n    @{Find = \"626\u00d73=\"; Replace = \"492\u00d73=\"},\n    @{Find = \"102\u00d78=\"; Replace = \"659\u00d77=\"},\n    @{Find = \"664\u00d76=\"; Replace = \"511\u00d78=\"},\n    @{Find = \"223\u00d73=\"; Replace = \"399\u00d74=\"},\n    @{Find = \"626\u00d75=\"; Replace = \"874\u00d72=\"},\n    @{Find = \"428\u00d72=\"; Replace = \"766\u00d78=\"},\n    @{Find = \"954\u00d77=\"; Replace = \"336\u00d72=\"},\n    @{Find = \"297\u00d72=\"; Replace = \"149\u00d75=\"},\n    @{Find = \"464\u00d77=\"; Replace = \"573\u00d74=\"},\n    @{Find = \"679\u00d75=\"; Replace = \"221\u00d76=\"},\n    @{Find = \"913\u00d75=\"; Replace = \"322\u00d74=\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
